$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5789666666666667
$ws.Range("H2").Value = 1.7369
$ws.Range("I2").Value = 0.01523705650035473
$ws.Range("J2").Value = 0.01523705650035472
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 6.453652272466669
$ws.Range("R2").Value = 58.08287045220001
$ws.Range("S2").Value = 0.003953720693485062
$ws.Range("T2").Value = 0.003953720693485061
$ws.Range("G3").Value = 0.5789666666666667
$ws.Range("H3").Value = 1.7369
$ws.Range("I3").Value = 0.01523705650035473
$ws.Range("J3").Value = 0.01523705650035472
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 16.0289799131
$ws.Range("R3").Value = 144.2608192179
$ws.Range("S3").Value = 0.009819882897666171
$ws.Range("T3").Value = 0.009819882897666169
$ws.Range("G4").Value = 0.5789666666666667
$ws.Range("H4").Value = 1.7369
$ws.Range("I4").Value = 0.01523705650035473
$ws.Range("J4").Value = 0.01523705650035472
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 2.388791957077778
$ws.Range("R4").Value = 21.4991276137
$ws.Range("S4").Value = 0.001463452909203495
$ws.Range("T4").Value = 0.001463452909203494
$ws.Range("I5").Value = 0.6545086962501954
$ws.Range("J5").Value = 0.6545086962501954
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 277.2170290768387
$ws.Range("R5").Value = 2494.953261691548
$ws.Range("S5").Value = 0.1698323148155342
$ws.Range("T5").Value = 0.1698323148155342
$ws.Range("I6").Value = 0.6545086962501954
$ws.Range("J6").Value = 0.6545086962501954
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.4218136719865446
$ws.Range("T6").Value = 0.4218136719865446
$ws.Range("I7").Value = 0.6545086962501954
$ws.Range("J7").Value = 0.6545086962501954
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.06286270944811656
$ws.Range("T7").Value = 0.06286270944811656
$ws.Range("I8").Value = 0.33025424724945
$ws.Range("J8").Value = 0.3302542472494499
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 139.8791212202674
$ws.Range("R8").Value = 1258.912090982406
$ws.Range("S8").Value = 0.08569457305819429
$ws.Range("T8").Value = 0.08569457305819428
$ws.Range("I9").Value = 0.33025424724945
$ws.Range("J9").Value = 0.3302542472494499
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.2128401922228869
$ws.Range("T9").Value = 0.2128401922228869
$ws.Range("I10").Value = 0.33025424724945
$ws.Range("J10").Value = 0.3302542472494499
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.03171948196836876
$ws.Range("T10").Value = 0.03171948196836875
